$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mobile number used for sendkeys changed
$ws.Range("A2").Value = "9988899999"

# File path fetch switched from a hardcoded absolute OS path to a
# relative "InputFiles\..." path
$ws.Range("S2").Value = "InputFiles\GoldLoan.jpeg"
$ws.Range("T2").Value = "InputFiles\GoldLoan.jpeg"
$ws.Range("Y2").Value = "InputFiles\GoldLoan.jpeg"
$ws.Range("AD2").Value = "InputFiles\GoldLoan.jpeg"

# Update the view: scrolled right so column K is the left-most visible
# column, and the active selection moved to AD2
$ws.Range("AD2").Select()
$excel.ActiveWindow.ScrollColumn = 11
